$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 140, shifting existing rows 140-252 down to 141-253
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly record
$ws.Cells.Item(140, 1).Value = 8
$ws.Cells.Item(140, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(140, 3).Value = "Coquimbo"
$ws.Cells.Item(140, 4).Value = 45096
$ws.Cells.Item(140, 5).Value = 4
$ws.Cells.Item(140, 6).Value = 100112001
$ws.Cells.Item(140, 7).Value = "Berenjena"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 400
$ws.Cells.Item(140, 11).Value = 7000
$ws.Cells.Item(140, 12).Value = 8000
$ws.Cells.Item(140, 13).Value = 7500
$ws.Cells.Item(140, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(140, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(140, 16).Value = 150
$ws.Cells.Item(140, 17).Value = 50
$ws.Cells.Item(140, 18).Value = "Hortaliza"
